$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add the new "alu" result block (rows 14-16), mirroring the existing
# "Hamming" block (rows 5-6): row 14 is a blank spacer row with only the
# bold styled U column, row 15 carries the "alu" label + B=4 sample data,
# row 16 carries the B=8 sample data. Columns P..U use the same SUM/weighted
# formulas as the rest of the sheet (library weight factor of 5 for S).
# ---------------------------------------------------------------------------

# Row 14 - spacer row (only U14 present, bold like the other spacer rows)
$ws.Range("U14").Font.Bold = $true

# Row 15 - "alu", B=4
$ws.Range("A15").Value = "alu"
$ws.Range("A15").Font.Bold = $true
$ws.Range("B15").Value = 4
$ws.Range("B15").Font.Bold = $true
$ws.Range("P15").Value = 26
$ws.Range("Q15").Value = 81
$ws.Range("R15").Formula = "=SUM(P15:Q15)"
$ws.Range("S15").Value = 163
$ws.Range("T15").Formula = "=SUM(P15:R15)"
$ws.Range("U15").Formula = "=SUM(P15:Q15)+5*S15"
$ws.Range("U15").Font.Bold = $true

# Row 16 - "alu", B=8
$ws.Range("B16").Value = 8
$ws.Range("B16").Font.Bold = $true
$ws.Range("P16").Value = 139
$ws.Range("Q16").Value = 220
$ws.Range("R16").Formula = "=SUM(P16:Q16)"
$ws.Range("S16").Value = 422
$ws.Range("T16").Formula = "=SUM(P16:R16)"
$ws.Range("U16").Formula = "=SUM(P16:Q16)+5*S16"
$ws.Range("U16").Font.Bold = $true

# Match the author's final cursor position recorded in the saved workbook.
$ws.Range("U20").Select() | Out-Null
